$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new data rows for the new weekly report date, right above the
# existing row 1290, shifting all rows from 1290 downward to 1295 onward.
$ws.Rows("1290:1294").Insert()

# New row 1290: Sandia, Calidad "Extra"
$ws.Range("A1290").Value = 6
$ws.Range("B1290").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1290").Value = "Metropolitana"
$ws.Range("D1290").Value = 44939
$ws.Range("E1290").Value = 13
$ws.Range("F1290").Value = 100112028
$ws.Range("G1290").Value = "Sandia"
$ws.Range("H1290").Value = "Sin especificar"
$ws.Range("I1290").Value = "Extra"
$ws.Range("J1290").Value = 2900
$ws.Range("K1290").Value = 3000
$ws.Range("L1290").Value = 3300
$ws.Range("M1290").Value = 3155
$ws.Range("N1290").Value = "$/unidad"
$ws.Range("O1290").Value = "Región de O'Higgins"
$ws.Range("P1290").Value = 3155
$ws.Range("Q1290").Value = 1
$ws.Range("R1290").Value = "Hortaliza"

# New row 1291: Sandia, Calidad "Primera"
$ws.Range("A1291").Value = 6
$ws.Range("B1291").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1291").Value = "Metropolitana"
$ws.Range("D1291").Value = 44939
$ws.Range("E1291").Value = 13
$ws.Range("F1291").Value = 100112028
$ws.Range("G1291").Value = "Sandia"
$ws.Range("H1291").Value = "Sin especificar"
$ws.Range("I1291").Value = "Primera"
$ws.Range("J1291").Value = 6400
$ws.Range("K1291").Value = 2500
$ws.Range("L1291").Value = 2700
$ws.Range("M1291").Value = 2591
$ws.Range("N1291").Value = "$/unidad"
$ws.Range("O1291").Value = "Región de O'Higgins"
$ws.Range("P1291").Value = 2591
$ws.Range("Q1291").Value = 1
$ws.Range("R1291").Value = "Hortaliza"

# New row 1292: Sandia, Calidad "Segunda"
$ws.Range("A1292").Value = 6
$ws.Range("B1292").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1292").Value = "Metropolitana"
$ws.Range("D1292").Value = 44939
$ws.Range("E1292").Value = 13
$ws.Range("F1292").Value = 100112028
$ws.Range("G1292").Value = "Sandia"
$ws.Range("H1292").Value = "Sin especificar"
$ws.Range("I1292").Value = "Segunda"
$ws.Range("J1292").Value = 5100
$ws.Range("K1292").Value = 1700
$ws.Range("L1292").Value = 2000
$ws.Range("M1292").Value = 1871
$ws.Range("N1292").Value = "$/unidad"
$ws.Range("O1292").Value = "Región de O'Higgins"
$ws.Range("P1292").Value = 1871
$ws.Range("Q1292").Value = 1
$ws.Range("R1292").Value = "Hortaliza"

# New row 1293: Sandia, Calidad "Super"
$ws.Range("A1293").Value = 6
$ws.Range("B1293").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1293").Value = "Metropolitana"
$ws.Range("D1293").Value = 44939
$ws.Range("E1293").Value = 13
$ws.Range("F1293").Value = 100112028
$ws.Range("G1293").Value = "Sandia"
$ws.Range("H1293").Value = "Sin especificar"
$ws.Range("I1293").Value = "Super"
$ws.Range("J1293").Value = 4200
$ws.Range("K1293").Value = 3500
$ws.Range("L1293").Value = 3800
$ws.Range("M1293").Value = 3636
$ws.Range("N1293").Value = "$/unidad"
$ws.Range("O1293").Value = "Región de O'Higgins"
$ws.Range("P1293").Value = 3636
$ws.Range("Q1293").Value = 1
$ws.Range("R1293").Value = "Hortaliza"

# New row 1294: Sandia, Calidad "Tercera"
$ws.Range("A1294").Value = 6
$ws.Range("B1294").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1294").Value = "Metropolitana"
$ws.Range("D1294").Value = 44939
$ws.Range("E1294").Value = 13
$ws.Range("F1294").Value = 100112028
$ws.Range("G1294").Value = "Sandia"
$ws.Range("H1294").Value = "Sin especificar"
$ws.Range("I1294").Value = "Tercera"
$ws.Range("J1294").Value = 3700
$ws.Range("K1294").Value = 1300
$ws.Range("L1294").Value = 1500
$ws.Range("M1294").Value = 1424
$ws.Range("N1294").Value = "$/unidad"
$ws.Range("O1294").Value = "Región de O'Higgins"
$ws.Range("P1294").Value = 1424
$ws.Range("Q1294").Value = 1
$ws.Range("R1294").Value = "Hortaliza"
